$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "id_expressao"

# Update the sheet's default (standard) row height: 15 -> 15.75
$ws.StandardHeight = 15.75

# Update selection to C4
$ws.Range("C4").Select()
